$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row range that holds data (header is row 1, data starts row 2..31)
$firstRow = 2
$lastRow = 31
$lastCol = 26   # column Z

# Mapping: target (after) row -> source (before) row.
# Derived from the unified diff: the whole data table got re-fetched / re-ordered,
# each record (all of its columns) moving intact to a new row position.
$rowMap = @{
    2 = 2;  3 = 3;
    4 = 7;  5 = 5;  6 = 4;  7 = 6;
    8 = 8;  9 = 9;  10 = 10; 11 = 11;
    12 = 25; 13 = 13; 14 = 20; 15 = 27; 16 = 23; 17 = 12; 18 = 24; 19 = 26;
    20 = 30; 21 = 29; 22 = 31; 23 = 16; 24 = 17; 25 = 21; 26 = 19; 27 = 28;
    28 = 22; 29 = 18; 30 = 15; 31 = 14
}

# New "Förändrad" (column C) date serial for every data row (45077 -> 46078,
# i.e. 2026-02-24 -> 2026-02-25).
$newChangedDate = 46078

# ---- Step 1: snapshot every cell (formula or literal) of every source row ----
# Captured up-front, before any writes, so overlapping source/target rows never
# read back already-overwritten data.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula()) {
            $rowData[$c] = @{ kind = "formula"; data = $cell.Formula }
        } else {
            $rowData[$c] = @{ kind = "value"; data = $cell.Value2() }
        }
    }
    $snapshot[$r] = $rowData
}

# ---- Step 2: write every target row from its mapped source row's snapshot ----
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $rowData = $snapshot[$sourceRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $entry = $rowData[$c]
        $destCell = $ws.Cells.Item($targetRow, $c)
        if ($entry.kind -eq "formula") {
            $destCell.Formula = $entry.data
        } else {
            $destCell.Value2 = $entry.data
        }
    }
}

# ---- Step 3: bump the "Förändrad" date for every data row ----
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value2 = $newChangedDate
}
